$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B11 held the "R40" label; the sheet now stores the text "1" there instead
# (a new shared string). Prefix with an apostrophe so the numeric-looking
# text is stored as text, not coerced into a number.
$ws.Range("B11").Value = "'1"
